$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Matriz_Resultados")
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("E6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("B10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0

$ws = $wb.Worksheets.Item("P_valores")
$ws.Range("C2").Value = 0.0005099866019993016
$ws.Range("D2").Value = 0.0006863075709535593
$ws.Range("E2").Value = 0.002954522280926408
$ws.Range("F2").Value = 0.0144398756629347
$ws.Range("G2").Value = 0.0003629518155701117
$ws.Range("H2").Value = 0.000624128788740963
$ws.Range("I2").Value = [double]"2.032366721205214E-05"
$ws.Range("J2").Value = 0.009378683866079207
$ws.Range("B3").Value = 0.0005099866019993016
$ws.Range("D3").Value = [double]"4.374867511147684E-06"
$ws.Range("E3").Value = 0.0003097350608742833
$ws.Range("F3").Value = 0.0004501522631856769
$ws.Range("G3").Value = 0.00424346945235432
$ws.Range("H3").Value = 0.001458143148378221
$ws.Range("I3").Value = 0.04841437440638074
$ws.Range("J3").Value = [double]"2.852206422154779E-07"
$ws.Range("B4").Value = 0.0006863075709535593
$ws.Range("C4").Value = [double]"4.374867511147684E-06"
$ws.Range("E4").Value = 0.0008614675191416232
$ws.Range("F4").Value = 0.0006227855471649946
$ws.Range("G4").Value = 0.01395811335666663
$ws.Range("H4").Value = 0.006765181803076281
$ws.Range("I4").Value = 0.0939954178601945
$ws.Range("J4").Value = [double]"3.024653429939406E-07"
$ws.Range("B5").Value = 0.002954522280926408
$ws.Range("C5").Value = 0.0003097350608742833
$ws.Range("D5").Value = 0.0008614675191416232
$ws.Range("F5").Value = 0.003201539792706987
$ws.Range("G5").Value = 0.7139669579277137
$ws.Range("H5").Value = 0.3708102026103561
$ws.Range("I5").Value = 0.8644818850851959
$ws.Range("J5").Value = 0.002399303045043544
$ws.Range("B6").Value = 0.0144398756629347
$ws.Range("C6").Value = 0.0004501522631856769
$ws.Range("D6").Value = 0.0006227855471649946
$ws.Range("E6").Value = 0.003201539792706987
$ws.Range("G6").Value = 0.0002899130442952647
$ws.Range("H6").Value = 0.0005705696499835877
$ws.Range("I6").Value = [double]"1.207479872555872E-05"
$ws.Range("J6").Value = 0.0119971898759561
$ws.Range("B7").Value = 0.0003629518155701117
$ws.Range("C7").Value = 0.00424346945235432
$ws.Range("D7").Value = 0.01395811335666663
$ws.Range("E7").Value = 0.7139669579277137
$ws.Range("F7").Value = 0.0002899130442952647
$ws.Range("H7").Value = 0.3370891896476373
$ws.Range("I7").Value = 0.5171735538696869
$ws.Range("J7").Value = 0.004748961281126673
$ws.Range("B8").Value = 0.000624128788740963
$ws.Range("C8").Value = 0.001458143148378221
$ws.Range("D8").Value = 0.006765181803076281
$ws.Range("E8").Value = 0.3708102026103561
$ws.Range("F8").Value = 0.0005705696499835877
$ws.Range("G8").Value = 0.3370891896476373
$ws.Range("I8").Value = 0.4288856157437786
$ws.Range("J8").Value = 0.0001936995406381925
$ws.Range("B9").Value = [double]"2.032366721205214E-05"
$ws.Range("C9").Value = 0.04841437440638074
$ws.Range("D9").Value = 0.0939954178601945
$ws.Range("E9").Value = 0.8644818850851959
$ws.Range("F9").Value = [double]"1.207479872555872E-05"
$ws.Range("G9").Value = 0.5171735538696869
$ws.Range("H9").Value = 0.4288856157437786
$ws.Range("J9").Value = 0.246855690628649
$ws.Range("B10").Value = 0.009378683866079207
$ws.Range("C10").Value = [double]"2.852206422154779E-07"
$ws.Range("D10").Value = [double]"3.024653429939406E-07"
$ws.Range("E10").Value = 0.002399303045043544
$ws.Range("F10").Value = 0.0119971898759561
$ws.Range("G10").Value = 0.004748961281126673
$ws.Range("H10").Value = 0.0001936995406381925
$ws.Range("I10").Value = 0.246855690628649

$ws = $wb.Worksheets.Item("Estadisticos_DM")
$ws.Range("C2").Value = 4.488825824202573
$ws.Range("D2").Value = 4.334435018164408
$ws.Range("E2").Value = 3.590349634616223
$ws.Range("F2").Value = 2.790782439492524
$ws.Range("G2").Value = 4.667340732493028
$ws.Range("H2").Value = 4.3836733686979
$ws.Range("I2").Value = 6.276995050029293
$ws.Range("J2").Value = 3.0092115351993
$ws.Range("B3").Value = -4.488825824202573
$ws.Range("D3").Value = -7.22735218614586
$ws.Range("E3").Value = -4.75122377276956
$ws.Range("F3").Value = -4.554111045667168
$ws.Range("G3").Value = -3.408157909227119
$ws.Range("H3").Value = -3.947766999416877
$ws.Range("I3").Value = -2.162081006760383
$ws.Range("J3").Value = -9.129182374361463
$ws.Range("B4").Value = -4.334435018164408
$ws.Range("C4").Value = 7.22735218614586
$ws.Range("E4").Value = -4.217080848319907
$ws.Range("F4").Value = -4.38479189519682
$ws.Range("G4").Value = -2.808023278408015
$ws.Range("H4").Value = -3.17370585706138
$ws.Range("I4").Value = -1.796634031314869
$ws.Range("J4").Value = -9.085007528599535
$ws.Range("B5").Value = -3.590349634616223
$ws.Range("C5").Value = 4.75122377276956
$ws.Range("D5").Value = 4.217080848319907
$ws.Range("F5").Value = -3.549903625197971
$ws.Range("G5").Value = 0.3740558391805323
$ws.Range("H5").Value = 0.9246534502568798
$ws.Range("I5").Value = -0.173837323908816
$ws.Range("J5").Value = -3.695351410791463
$ws.Range("B6").Value = -2.790782439492524
$ws.Range("C6").Value = 4.554111045667168
$ws.Range("D6").Value = 4.38479189519682
$ws.Range("E6").Value = 3.549903625197971
$ws.Range("G6").Value = 4.786344920236681
$ws.Range("H6").Value = 4.430309985567148
$ws.Range("I6").Value = 6.590779123773242
$ws.Range("J6").Value = 2.884789777860116
$ws.Range("B7").Value = -4.667340732493028
$ws.Range("C7").Value = 3.408157909227119
$ws.Range("D7").Value = 2.808023278408015
$ws.Range("E7").Value = -0.3740558391805323
$ws.Range("F7").Value = -4.786344920236681
$ws.Range("H7").Value = 0.9940087980018663
$ws.Range("I7").Value = -0.664488758506576
$ws.Range("J7").Value = -3.351585467001053
$ws.Range("B8").Value = -4.3836733686979
$ws.Range("C8").Value = 3.947766999416877
$ws.Range("D8").Value = 3.17370585706138
$ws.Range("E8").Value = -0.9246534502568798
$ws.Range("F8").Value = -4.430309985567148
$ws.Range("G8").Value = -0.9940087980018663
$ws.Range("I8").Value = -0.8147000302543833
$ws.Range("J8").Value = -5.00226644083312
$ws.Range("B9").Value = -6.276995050029293
$ws.Range("C9").Value = 2.162081006760383
$ws.Range("D9").Value = 1.796634031314869
$ws.Range("E9").Value = 0.173837323908816
$ws.Range("F9").Value = -6.590779123773242
$ws.Range("G9").Value = 0.664488758506576
$ws.Range("H9").Value = 0.8147000302543833
$ws.Range("J9").Value = -1.208543195606404
$ws.Range("B10").Value = -3.0092115351993
$ws.Range("C10").Value = 9.129182374361463
$ws.Range("D10").Value = 9.085007528599535
$ws.Range("E10").Value = 3.695351410791463
$ws.Range("F10").Value = -2.884789777860116
$ws.Range("G10").Value = 3.351585467001053
$ws.Range("H10").Value = 5.00226644083312
$ws.Range("I10").Value = 1.208543195606404

$ws = $wb.Worksheets.Item("Resumen")
$ws.Range("B2").Value = 5
$ws.Range("D2").Value = 3
$ws.Range("E2").Value = 62.5
$ws.Range("B3").Value = 4
$ws.Range("D3").Value = 3
$ws.Range("E3").Value = 50
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 5
$ws.Range("A5").Value = "DeepAR"
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 6
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 2.721317932732118
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 6
$ws.Range("E6").Value = 25
$ws.Range("A7").Value = "Block Bootstrapping"
$ws.Range("B7").Value = 0
$ws.Range("D7").Value = 3
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 8.793187674308667
$ws.Range("A8").Value = "LSPMW"
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 2
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 2.552051746246231
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 5
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 0
$ws.Range("A10").Value = "EnCQR-LSTM"
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 5
$ws.Range("F10").Value = 3.838167971691867

Write-Output "edit applied"
